$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "SAP number" label row is being renamed to "Vendor number".
$ws.Range("A3").Value = "Vendor number"

# Leave the cursor/selection on the cell that was edited, matching what
# Excel records as the active cell when the workbook is saved.
$ws.Range("A3").Select()
